{"js": "// Find the paragraph that holds the (partially split) Global Fishing Watch\n// citation + hyperlink. Locating it via search (instead of a hard-coded\n// paragraph index) keeps this robust to unrelated structural differences.\nconst searchResults = context.document.body.search(\"Global Fishing Watch\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the 'Global Fishing Watch' citation paragraph.\");\n}\n\nconst citationPara = searchResults.items[0].paragraphs.getFirst();\ncitationPara.load(\"text\");\nawait context.sync();\n\nconst fishingWatchUrl =\n  \"https://globalfishingwatch.org/data-download/datasets/public-fishing-vessels-v1\";\n\n// The hyperlink text in the original document is split across three runs\n// (\"https://globalfishi\" + \"n\" + \"gwatch.org/...\").  Office.js has no\n// \"merge runs\" primitive, so rebuild the whole paragraph's content: the\n// plain lead-in text, then the hyperlink as a single run, then the\n// trailing \". \" text \u2014 preserving the paragraph's existing text exactly.\ncitationPara.clear();\nawait context.sync();\n\ncitationPara.insertText(\"Global Fishing Watch. 2020. Accessed on July 22, 2020. \", Word.InsertLocation.end);\nawait context.sync();\n\nconst newHyperlinkRange = citationPara.insertText(fishingWatchUrl, Word.InsertLocation.end);\nawait context.sync();\nnewHyperlinkRange.hyperlink = fishingWatchUrl;\nawait context.sync();\n\ncitationPara.insertText(\". \", Word.InsertLocation.end);\nawait context.sync();\n\n// Add a blank separator paragraph after the citation.\nconst blankPara = citationPara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// Add the \"Annual fishing totals:\" heading line.\nconst totalsHeadingPara = blankPara.insertParagraph(\"Annual fishing totals:\", Word.InsertLocation.after);\nawait context.sync();\n\n// Add the new FAO hyperlink paragraph.\nconst faoPara = totalsHeadingPara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst faoUrl = \"http://www.fao.org/fishery/statistics/global-capture-production/en\";\nconst faoHyperlinkRange = faoPara.insertText(faoUrl, Word.InsertLocation.start);\nawait context.sync();\nfaoHyperlinkRange.hyperlink = faoUrl;\nawait context.sync();\n\n// Trailing blank paragraph, before the final section break.\nfaoPara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1. Merge the split \"https://globalfishi\" + \"n\" + \"gwatch.org/...\" runs of\n#    the existing Global Fishing Watch hyperlink into a single run. Locate\n#    the hyperlink by its target address (robust to index changes) rather\n#    than assuming it is the last one, then setting TextToDisplay rewrites\n#    the hyperlink's run(s) as one run while leaving the surrounding text\n#    (\"... Accessed on ... . \" and \". \") untouched.\n# ---------------------------------------------------------------------------\n$fishingWatchUrl = \"https://globalfishingwatch.org/data-download/datasets/public-fishing-vessels-v1\"\n$fishingWatchLink = $null\nfor ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {\n    $candidate = $d.Hyperlinks.Item($i)\n    if ($candidate.Address -eq $fishingWatchUrl) {\n        $fishingWatchLink = $candidate\n    }\n}\nif ($fishingWatchLink -eq $null) {\n    throw \"Could not find the Global Fishing Watch hyperlink.\"\n}\n$fishingWatchLink.TextToDisplay = $fishingWatchUrl\n\n# ---------------------------------------------------------------------------\n# 2. Append, after the citation paragraph: a blank paragraph, an\n#    \"Annual fishing totals:\" paragraph, a paragraph with a new hyperlink to\n#    the FAO capture-production statistics, and a trailing blank paragraph.\n# ---------------------------------------------------------------------------\n\n# Blank separator paragraph.\n$end = $d.Content.End\n$r = $d.Range($end, $end)\n$r.InsertParagraphAfter()\n\n# \"Annual fishing totals:\" paragraph.\n$end = $d.Content.End\n$r = $d.Range($end, $end)\n$r.InsertParagraphAfter()\n$insPoint = $d.Content.End - 1\n$r2 = $d.Range($insPoint, $insPoint)\n$r2.InsertAfter(\"Annual fishing totals:\")\n\n# Paragraph with the new FAO hyperlink.\n$end = $d.Content.End\n$r = $d.Range($end, $end)\n$r.InsertParagraphAfter()\n$insPoint = $d.Content.End - 1\n$faoUrl = \"http://www.fao.org/fishery/statistics/global-capture-production/en\"\n$r2 = $d.Range($insPoint, $insPoint)\n$r2.InsertAfter($faoUrl)\n$linkRange = $d.Range($insPoint, $insPoint + $faoUrl.Length)\n$d.Hyperlinks.Add($linkRange, $faoUrl) | Out-Null\n\n# Trailing blank paragraph.\n$end = $d.Content.End\n$r = $d.Range($end, $end)\n$r.InsertParagraphAfter()\n"}
